# Apply updated Leve profit-calculation figures across the Sheets workbook.
# Values below were computed by the scheduled profits runner and mirror the
# authoritative OOXML diff for this commit.
$wb = $excel.ActiveWorkbook

# ==== Sheet: ALC ====
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1247.7894
$ws.Range("I6").Value = 231.75
$ws.Range("J6").Value = 6666.6665
$ws.Range("K6").Value = 695.25
$ws.Range("L6").Value = 19999.9995
$ws.Range("M6").Value = -583.25
$ws.Range("N6").Value = -20223.9995

# Row 44
$ws.Range("H44").Value = 15666.667
$ws.Range("J44").Value = 15666.667
$ws.Range("L44").Value = 15666.667
$ws.Range("N44").Value = -16590.667

# Row 132
$ws.Range("H132").Value = 2377.8408
$ws.Range("I132").Value = 1815.1923
$ws.Range("K132").Value = 5445.5769
$ws.Range("M132").Value = -2915.5769

# Row 137
$ws.Range("H137").Value = 1018.89655
$ws.Range("J137").Value = 1337.7727
$ws.Range("L137").Value = 4013.3181
$ws.Range("N137").Value = -9113.3181

# Row 138
$ws.Range("H138").Value = 3694.2222
$ws.Range("I138").Value = 1777.92
$ws.Range("J138").Value = 8049.4546
$ws.Range("K138").Value = 5333.76
$ws.Range("L138").Value = 24148.3638
$ws.Range("M138").Value = -193.7600000000002
$ws.Range("N138").Value = -34428.3638

# Row 141
$ws.Range("H141").Value = 3174.9524
$ws.Range("I141").Value = 1530.7966
$ws.Range("K141").Value = 4592.3898
$ws.Range("M141").Value = 587.6102000000001

# ==== Sheet: ARM ====
$ws = $wb.Worksheets.Item("ARM")
# Row 12
$ws.Range("H12").Value = 10000
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

# Row 61
$ws.Range("H61").Value = 1422.8889
$ws.Range("I61").Value = 1149.6
$ws.Range("J61").Value = 2789.3333
$ws.Range("K61").Value = 1149.6
$ws.Range("L61").Value = 2789.3333
$ws.Range("M61").Value = -937.5999999999999
$ws.Range("N61").Value = -3213.3333

# Row 132
$ws.Range("H132").Value = 1881.2927
$ws.Range("I132").Value = 1132.1154
$ws.Range("K132").Value = 3396.3462
$ws.Range("M132").Value = -866.3462

# Row 136
$ws.Range("H136").Value = 1422.8889
$ws.Range("I136").Value = 1149.6
$ws.Range("J136").Value = 2789.3333
$ws.Range("K136").Value = 3448.8
$ws.Range("L136").Value = 8367.999899999999
$ws.Range("M136").Value = -898.7999999999997
$ws.Range("N136").Value = -13467.9999

# ==== Sheet: BSM ====
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1825.1277
$ws.Range("I134").Value = 1550.9756
$ws.Range("J134").Value = 3698.5
$ws.Range("K134").Value = 4652.9268
$ws.Range("L134").Value = 11095.5
$ws.Range("M134").Value = -2117.9268
$ws.Range("N134").Value = -16165.5

# ==== Sheet: CRP ====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1747.68
$ws.Range("I31").Value = 1233.579
$ws.Range("K31").Value = 1233.579
$ws.Range("M31").Value = -938.579

# Row 34
$ws.Range("H34").Value = 1747.68
$ws.Range("I34").Value = 1233.579
$ws.Range("K34").Value = 1233.579
$ws.Range("M34").Value = -1031.579

# Row 58
$ws.Range("H58").Value = 741937.6
$ws.Range("I58").Value = 1372851.2
$ws.Range("J58").Value = 1299.8695
$ws.Range("K58").Value = 1372851.2
$ws.Range("L58").Value = 1299.8695
$ws.Range("M58").Value = -1372648.2
$ws.Range("N58").Value = -1705.8695

# Row 99
$ws.Range("H99").Value = 4008.6667
$ws.Range("I99").Value = 5214.2856
$ws.Range("J99").Value = 2320.8
$ws.Range("K99").Value = 5214.2856
$ws.Range("L99").Value = 2320.8
$ws.Range("M99").Value = -3716.2856
$ws.Range("N99").Value = -5316.8

# Row 126
$ws.Range("H126").Value = 4008.6667
$ws.Range("I126").Value = 5214.2856
$ws.Range("J126").Value = 2320.8
$ws.Range("K126").Value = 15642.8568
$ws.Range("L126").Value = 6962.400000000001
$ws.Range("M126").Value = -13172.8568
$ws.Range("N126").Value = -11902.4

# Row 132
$ws.Range("H132").Value = 251614.06
$ws.Range("I132").Value = 322677.4
$ws.Range("K132").Value = 968032.2000000001
$ws.Range("M132").Value = -965502.2000000001

# Row 134
$ws.Range("H134").Value = 1759.381
$ws.Range("I134").Value = 1356.3636
$ws.Range("J134").Value = 3237.111
$ws.Range("K134").Value = 4069.0908
$ws.Range("L134").Value = 9711.332999999999
$ws.Range("M134").Value = -1534.0908
$ws.Range("N134").Value = -14781.333

# Row 136
$ws.Range("H136").Value = 741937.6
$ws.Range("I136").Value = 1372851.2
$ws.Range("J136").Value = 1299.8695
$ws.Range("K136").Value = 4118553.6
$ws.Range("L136").Value = 3899.6085
$ws.Range("M136").Value = -4116003.6
$ws.Range("N136").Value = -8999.6085

# ==== Sheet: CUL ====
$ws = $wb.Worksheets.Item("CUL")
# Row 51
$ws.Range("H51").Value = 1300
$ws.Range("I51").Value = 900
$ws.Range("J51").Value = 1500
$ws.Range("K51").Value = 2700
$ws.Range("L51").Value = 4500
$ws.Range("M51").Value = -2240
$ws.Range("N51").Value = -5420

# ==== Sheet: GSM ====
$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 5909.143
$ws.Range("I9").Value = 260.66666
$ws.Range("J9").Value = 39800
$ws.Range("K9").Value = 260.66666
$ws.Range("L9").Value = 39800
$ws.Range("M9").Value = -90.66665999999998
$ws.Range("N9").Value = -40140

# Row 123
$ws.Range("H123").Value = 11784.857
$ws.Range("J123").Value = 11784.857
$ws.Range("L123").Value = 11784.857
$ws.Range("N123").Value = -16684.857

# Row 124
$ws.Range("H124").Value = 69780
$ws.Range("J124").Value = 69780
$ws.Range("L124").Value = 69780
$ws.Range("N124").Value = -79600

# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# Row 131
$ws.Range("H131").Value = 38277.2
$ws.Range("J131").Value = 38277.2
$ws.Range("L131").Value = 38277.2
$ws.Range("N131").Value = -48357.2

# Row 132
$ws.Range("H132").Value = 1185.4857
$ws.Range("I132").Value = 695.3077
$ws.Range("K132").Value = 2085.9231
$ws.Range("M132").Value = 444.0769

# ==== Sheet: LTW ====
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 3161.102
$ws.Range("I132").Value = 2766.257
$ws.Range("J132").Value = 4148.2144
$ws.Range("K132").Value = 8298.771000000001
$ws.Range("L132").Value = 12444.6432
$ws.Range("M132").Value = -5768.771000000001
$ws.Range("N132").Value = -17504.6432

# Row 136
$ws.Range("H136").Value = 2831.419
$ws.Range("I136").Value = 2973.7368
$ws.Range("J136").Value = 2354.2354
$ws.Range("K136").Value = 8921.2104
$ws.Range("L136").Value = 7062.706200000001
$ws.Range("M136").Value = -6371.2104
$ws.Range("N136").Value = -12162.7062

# ==== Sheet: WVR ====
$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 28005.6
$ws.Range("J15").Value = 28005.6
$ws.Range("L15").Value = 28005.6
$ws.Range("N15").Value = -28581.6

# Row 132
$ws.Range("H132").Value = 855.7627
$ws.Range("I132").Value = 606
$ws.Range("J132").Value = 1658.5714
$ws.Range("K132").Value = 1818
$ws.Range("L132").Value = 4975.7142
$ws.Range("M132").Value = 712
$ws.Range("N132").Value = -10035.7142

# Row 136
$ws.Range("H136").Value = 1490.2559
$ws.Range("I136").Value = 1252.6666
$ws.Range("K136").Value = 3757.9998
$ws.Range("M136").Value = -1207.9998
